$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text replacements (non-numeric-looking strings) ---
$ws.Range("A1").Value = 'Scenario: Test Suite'
$ws.Range("F2").Value = 'gag32'
$ws.Range("E3").Value = '// 👉 Enable drag-and-drop row reordering
tableView.setRowFactory(tv -> {
    TableRow<TestStep> row = new TableRow<>();
    row.setOnDragDetected(event -> {
        if (!row.isEmpty()) {
            Dragboard db = row.startDragAndDrop(TransferMode.MOVE);
            ClipboardContent cc = new ClipboardContent();
            cc.putString(Integer.toString(row.getIndex()));
            db.setContent(cc);
            event.consume();
        }
    });
    row.setOnDragOver(event -> {
        Dragboard db = event.getDragboard();
        if (db.hasString()) {
            int draggedIndex = Integer.parseInt(db.getString());
            if (row.getIndex() != draggedIndex) {
                event.acceptTransferModes(TransferMode.MOVE);
                row.setStyle("-fx-background-color: lightgreen;");
            }
        }
        event.consume();
    });
    row.setOnDragExited(event -> row.setStyle(""));
    row.setOnDragDropped(event -> {
        Dragboard db = event.getDragboard();
        if (db.hasString()) {
            int draggedIndex = Integer.parseInt(db.getString());
            TestStep draggedStep = tableView.getItems().remove(draggedIndex);
            int dropIndex = row.isEmpty() ? tableView.getItems().size() : row.getIndex();
            tableView.getItems().add(dropIndex, draggedStep);
            tableView.getSelectionModel().select(dropIndex);
            event.setDropCompleted(true);
        }
        event.consume();
    });
    return row;
});
'
$ws.Range("B4").Value = 'WebServer'
$ws.Range("C4").Value = 'sendRequest'
$ws.Range("B5").Value = 'Database'
$ws.Range("C5").Value = 'executeQuery'

# Recompute row heights for rows whose content changed (e.g. multi-line text)
# so no stray custom row-height survives in the saved file.
$ws.Rows(1).AutoFit()
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
$ws.Rows(4).AutoFit()
$ws.Rows(5).AutoFit()

# --- Numeric-looking text values: must stay text (shared string), not become a number. ---
# Use a helper cell (H1, outside the used range) with a formula producing a text result,
# then paste-special (values only) into the destination so it keeps the default style (s=0)
# while still being stored as text. The helper cell is then removed (shift left) so the
# worksheet dimension / layout is restored.
$ws.Range("H1").Formula = '="3"'
$ws.Range("H1").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("H1").Delete(-4159)
$ws.Range("H1").Formula = '="3"'
$ws.Range("H1").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("H1").Delete(-4159)
$ws.Range("H1").Formula = '="4"'
$ws.Range("H1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("H1").Delete(-4159)

# --- Cells that become empty: keep the cell present (blank) without affecting styles ---
$ws.Range("F3").Formula = '=""'
$ws.Range("E4").Formula = '=""'
$ws.Range("F4").Formula = '=""'
$ws.Range("E5").Formula = '=""'

$excel.CutCopyMode = 0
